# TC07 Bento multi-filter test case: widen the pr_status filter in every
# Neo4j / stat query on the sheet so it also matches "Not Reported", in
# addition to the existing "Positive" value.
#   d.pr_status In ["Positive"]  ->  d.pr_status In ["Positive","Not Reported"]

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldClause = 'd.pr_status In ["Positive"]'
$newClause = 'd.pr_status In ["Positive","Not Reported"]'

# Column B holds the per-tab query (CasesTab / SamplesTab / FilesTab rows),
# column C holds the shared summary/stat query used by all three rows.
$cellsToUpdate = @("B2", "C2", "B3", "C3", "B4", "C4")

foreach ($addr in $cellsToUpdate) {
    $cell = $ws.Range($addr)
    $text = [string]$cell.Text
    $cell.Value = $text.Replace($oldClause, $newClause)
}
